$wb = $excel.ActiveWorkbook

# Helper: write a value that LOOKS like a number but must be stored as TEXT
# (matches the original file where these cells are shared-string / text cells,
# not numeric cells). We temporarily force a Text number format so Excel does
# not auto-convert the numeric-looking string into a real number, then restore
# the cell style back to "Normal" so no stray formatting is left behind.
function Set-TextValue($range, [string]$val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# NOTE: worksheet names are matched case-insensitively by Worksheets.Item(),
# and this workbook has both "Vector_bf" and "Vector_BF" sheets which would
# otherwise collide. Use the (1-based) sheet position instead, which matches
# the order declared in workbook.xml:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# ---------------------------------------------------------------------
# Sheet: Restricciones_del_follower
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Row 2 (J_0_L0_v)
$ws3.Range("A2").Value = "-1.7806039396874054 - 2x_1 + 1.9677837616607594y_1 + 1.5150401184535802y_2"
Set-TextValue $ws3.Range("B2") "4.280603939687405"
Set-TextValue $ws3.Range("D2") "0.79"
Set-TextValue $ws3.Range("E2") "0"
Set-TextValue $ws3.Range("F2") "2.2"

# Row 3 (J_0_LP_v)
$ws3.Range("A3").Value = "1.1496672299905422 + x_1 - 3x_2 - 0.3477103695949638y_1 + 0.09638334121977188y_2"
Set-TextValue $ws3.Range("B3") "-3.149667229990542"
Set-TextValue $ws3.Range("D3") "0.09"
Set-TextValue $ws3.Range("E3") "2.4"
Set-TextValue $ws3.Range("F3") "0"

# Row 4 (J_Ne_L0_v)
$ws3.Range("A4").Value = "-51.80611610564271 + x_1 + x_2 + 3.184817181354109y_1 + 8.276583363312918y_2"
Set-TextValue $ws3.Range("B4") "49.79611610564271"
Set-TextValue $ws3.Range("D4") "0.54"
Set-TextValue $ws3.Range("E4") "0"
Set-TextValue $ws3.Range("F4") "9.200000000000001"

# ---------------------------------------------------------------------
# Sheet: Punto_modificado
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "6.65"
Set-TextValue $ws4.Range("B2") "2.1"
Set-TextValue $ws4.Range("C2") "5.2"
Set-TextValue $ws4.Range("D2") "3.2"

# ---------------------------------------------------------------------
# Sheet: Vector_bf
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "0.7569434836203279"
Set-TextValue $ws5.Range("A3") "-6.674911210477084"

# ---------------------------------------------------------------------
# Sheet: Vector_BF
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-0.3999999999999999"
Set-TextValue $ws6.Range("A3") "6.199999999999999"
Set-TextValue $ws6.Range("A4") "0.33450488702791314"
Set-TextValue $ws6.Range("A5") "-0.2313200189274525"

# ---------------------------------------------------------------------
# Sheet: Vector_Alpha (these are real numeric cells, not text)
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 0.3725647919618656
$ws7.Range("A3").Value = 0.968207398829925
